$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-08-15 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-16 Friday", 2) | Out-Null

# Update the answer table, cell by cell (row, col are 1-based)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "95-17=78"
$t.Cell(1,2).Range.Text = "25+38=63"
$t.Cell(1,3).Range.Text = "38+6=44"
$t.Cell(1,4).Range.Text = "80-49=31"
$t.Cell(1,5).Range.Text = "87-49=38"
$t.Cell(2,1).Range.Text = "44+37=81"
$t.Cell(2,2).Range.Text = "58+8=66"
$t.Cell(2,3).Range.Text = "90-24=66"
$t.Cell(2,4).Range.Text = "28+33=61"
$t.Cell(2,5).Range.Text = "27+4=31"
$t.Cell(3,1).Range.Text = "83-65=18"
$t.Cell(3,2).Range.Text = "35+47=82"
$t.Cell(3,3).Range.Text = "23-8=15"
$t.Cell(3,4).Range.Text = "83-54=29"
$t.Cell(3,5).Range.Text = "82-15=67"
$t.Cell(4,1).Range.Text = "83-57=26"
$t.Cell(4,2).Range.Text = "48+19=67"
$t.Cell(4,3).Range.Text = "90-88=2"
$t.Cell(4,4).Range.Text = "82-23=59"
$t.Cell(4,5).Range.Text = "55+29=84"
$t.Cell(5,1).Range.Text = "20-16=4"
$t.Cell(5,2).Range.Text = "36+25=61"
$t.Cell(5,3).Range.Text = "78+13=91"
$t.Cell(5,4).Range.Text = "94-5=89"
$t.Cell(5,5).Range.Text = "82-63=19"
$t.Cell(6,1).Range.Text = "36+19=55"
$t.Cell(6,2).Range.Text = "13+78=91"
$t.Cell(6,3).Range.Text = "58+6=64"
$t.Cell(6,4).Range.Text = "37+36=73"
$t.Cell(6,5).Range.Text = "83-59=24"
$t.Cell(7,1).Range.Text = "5+78=83"
$t.Cell(7,2).Range.Text = "65-17=48"
$t.Cell(7,3).Range.Text = "70-9=61"
$t.Cell(7,4).Range.Text = "73-16=57"
$t.Cell(7,5).Range.Text = "62-25=37"
$t.Cell(8,1).Range.Text = "9+59=68"
$t.Cell(8,2).Range.Text = "81-72=9"
$t.Cell(8,3).Range.Text = "29+18=47"
$t.Cell(8,4).Range.Text = "28+23=51"
$t.Cell(8,5).Range.Text = "49+12=61"
$t.Cell(9,1).Range.Text = "37+59=96"
$t.Cell(9,2).Range.Text = "39+8=47"
$t.Cell(9,3).Range.Text = "90-89=1"
$t.Cell(9,4).Range.Text = "82-65=17"
$t.Cell(9,5).Range.Text = "52-39=13"
$t.Cell(10,1).Range.Text = "24-16=8"
$t.Cell(10,2).Range.Text = "80-51=29"
$t.Cell(10,3).Range.Text = "17+4=21"
$t.Cell(10,4).Range.Text = "28+3=31"
$t.Cell(10,5).Range.Text = "69+12=81"
$t.Cell(11,1).Range.Text = "65+29=94"
$t.Cell(11,2).Range.Text = "85-27=58"
$t.Cell(11,3).Range.Text = "34-28=6"
$t.Cell(11,4).Range.Text = "16+37=53"
$t.Cell(11,5).Range.Text = "58+34=92"
$t.Cell(12,1).Range.Text = "55+37=92"
$t.Cell(12,2).Range.Text = "4+38=42"
$t.Cell(12,3).Range.Text = "52-25=27"
$t.Cell(12,4).Range.Text = "64+17=81"
$t.Cell(12,5).Range.Text = "81-52=29"
$t.Cell(13,1).Range.Text = "57+9=66"
$t.Cell(13,2).Range.Text = "60-54=6"
$t.Cell(13,3).Range.Text = "61-6=55"
$t.Cell(13,4).Range.Text = "48+8=56"
$t.Cell(13,5).Range.Text = "5+77=82"
$t.Cell(14,1).Range.Text = "26+6=32"
$t.Cell(14,2).Range.Text = "48-39=9"
$t.Cell(14,3).Range.Text = "27+49=76"
$t.Cell(14,4).Range.Text = "41-3=38"
$t.Cell(14,5).Range.Text = "98-69=29"
$t.Cell(15,1).Range.Text = "5+29=34"
$t.Cell(15,2).Range.Text = "28+8=36"
$t.Cell(15,3).Range.Text = "68-49=19"
$t.Cell(15,4).Range.Text = "18+18=36"
$t.Cell(15,5).Range.Text = "78-19=59"
$t.Cell(16,1).Range.Text = "20-4=16"
$t.Cell(16,2).Range.Text = "29+36=65"
$t.Cell(16,3).Range.Text = "91-62=29"
$t.Cell(16,4).Range.Text = "82-8=74"
$t.Cell(16,5).Range.Text = "12-4=8"
$t.Cell(17,1).Range.Text = "85-36=49"
$t.Cell(17,2).Range.Text = "30-5=25"
$t.Cell(17,3).Range.Text = "27+65=92"
$t.Cell(17,4).Range.Text = "18+25=43"
$t.Cell(17,5).Range.Text = "35-19=16"
$t.Cell(18,1).Range.Text = "86-18=68"
$t.Cell(18,2).Range.Text = "28+47=75"
$t.Cell(18,3).Range.Text = "86+5=91"
$t.Cell(18,4).Range.Text = "49+19=68"
$t.Cell(18,5).Range.Text = "56-47=9"
$t.Cell(19,1).Range.Text = "67-58=9"
$t.Cell(19,2).Range.Text = "51-8=43"
$t.Cell(19,3).Range.Text = "78+19=97"
$t.Cell(19,4).Range.Text = "57-8=49"
$t.Cell(19,5).Range.Text = "54+9=63"
$t.Cell(20,1).Range.Text = "62-39=23"
$t.Cell(20,2).Range.Text = "63-19=44"
$t.Cell(20,3).Range.Text = "87-48=39"
$t.Cell(20,4).Range.Text = "66-57=9"
$t.Cell(20,5).Range.Text = "85-39=46"
